$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1203.3889
$ws.Range("I62").Value = 939.9231
$ws.Range("K62").Value = 939.9231
$ws.Range("M62").Value = -315.9231
$ws.Range("H65").Value = 1203.3889
$ws.Range("I65").Value = 939.9231
$ws.Range("K65").Value = 4699.6155
$ws.Range("M65").Value = -1579.6155
$ws.Range("H111").Value = 1598.5883
$ws.Range("I111").Value = 1806
$ws.Range("J111").Value = 1365.25
$ws.Range("K111").Value = 5418
$ws.Range("L111").Value = 4095.75
$ws.Range("M111").Value = -2351
$ws.Range("N111").Value = -10229.75
$ws.Range("H112").Value = 76924296
$ws.Range("J112").Value = 76924296
$ws.Range("L112").Value = 230772888
$ws.Range("N112").Value = -230775104
$ws.Range("H129").Value = 1031.88
$ws.Range("I129").Value = 383.55554
$ws.Range("K129").Value = 1150.66662
$ws.Range("M129").Value = 3849.33338
$ws.Range("H137").Value = 3497.86
$ws.Range("I137").Value = 4173.8125
$ws.Range("J137").Value = 2296.1667
$ws.Range("K137").Value = 12521.4375
$ws.Range("L137").Value = 6888.500100000001
$ws.Range("M137").Value = -9971.4375
$ws.Range("N137").Value = -11988.5001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7664.5317
$ws.Range("I32").Value = 4290.1387
$ws.Range("J32").Value = 23331.357
$ws.Range("K32").Value = 4290.1387
$ws.Range("L32").Value = 23331.357
$ws.Range("M32").Value = -4003.1387
$ws.Range("N32").Value = -23905.357
$ws.Range("H61").Value = 1904.4762
$ws.Range("I61").Value = 1518.6171
$ws.Range("J61").Value = 3037.9375
$ws.Range("K61").Value = 1518.6171
$ws.Range("L61").Value = 3037.9375
$ws.Range("M61").Value = -1306.6171
$ws.Range("N61").Value = -3461.9375
$ws.Range("H132").Value = 8334833.5
$ws.Range("I132").Value = 11199152
$ws.Range("J132").Value = 2271.818
$ws.Range("K132").Value = 33597456
$ws.Range("L132").Value = 6815.454000000001
$ws.Range("M132").Value = -33594926
$ws.Range("N132").Value = -11875.454
$ws.Range("H136").Value = 1904.4762
$ws.Range("I136").Value = 1518.6171
$ws.Range("J136").Value = 3037.9375
$ws.Range("K136").Value = 4555.8513
$ws.Range("L136").Value = 9113.8125
$ws.Range("M136").Value = -2005.8513
$ws.Range("N136").Value = -14213.8125

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 26220.023
$ws.Range("I134").Value = 29091.82
$ws.Range("J134").Value = 3820
$ws.Range("K134").Value = 87275.45999999999
$ws.Range("L134").Value = 11460
$ws.Range("M134").Value = -84740.45999999999
$ws.Range("N134").Value = -16530

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4275418
$ws.Range("I31").Value = 1610.7273
$ws.Range("J31").Value = 9806228
$ws.Range("K31").Value = 1610.7273
$ws.Range("L31").Value = 9806228
$ws.Range("M31").Value = -1315.7273
$ws.Range("N31").Value = -9806818
$ws.Range("H34").Value = 4275418
$ws.Range("I34").Value = 1610.7273
$ws.Range("J34").Value = 9806228
$ws.Range("K34").Value = 1610.7273
$ws.Range("L34").Value = 9806228
$ws.Range("M34").Value = -1408.7273
$ws.Range("N34").Value = -9806632
$ws.Range("H60").Value = 10598.4
$ws.Range("I60").Value = 5146.5
$ws.Range("J60").Value = 14233
$ws.Range("K60").Value = 5146.5
$ws.Range("L60").Value = 14233
$ws.Range("M60").Value = -4635.5
$ws.Range("N60").Value = -15255
$ws.Range("H134").Value = 2148.2354
$ws.Range("I134").Value = 2230.125
$ws.Range("J134").Value = 1951.7
$ws.Range("K134").Value = 6690.375
$ws.Range("L134").Value = 5855.1
$ws.Range("M134").Value = -4155.375
$ws.Range("N134").Value = -10925.1

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 709.9
$ws.Range("I4").Value = 549.625
$ws.Range("K4").Value = 1648.875
$ws.Range("M4").Value = -1536.875
$ws.Range("H68").Value = 858.96295
$ws.Range("J68").Value = 1176.8
$ws.Range("L68").Value = 3530.4
$ws.Range("N68").Value = -5152.4
$ws.Range("H71").Value = 858.96295
$ws.Range("J71").Value = 1176.8
$ws.Range("L71").Value = 10591.2
$ws.Range("N71").Value = -18703.2
$ws.Range("H107").Value = 24323.365
$ws.Range("I107").Value = 21824.021
$ws.Range("J107").Value = 27414.658
$ws.Range("K107").Value = 65472.063
$ws.Range("L107").Value = 82243.974
$ws.Range("M107").Value = -63552.063
$ws.Range("N107").Value = -86083.974
$ws.Range("H113").Value = 459.9697
$ws.Range("I113").Value = 499.44446
$ws.Range("J113").Value = 445.16666
$ws.Range("K113").Value = 1498.33338
$ws.Range("L113").Value = 1335.49998
$ws.Range("M113").Value = 671.66662
$ws.Range("N113").Value = -5675.499980000001
$ws.Range("H131").Value = 840.48
$ws.Range("I131").Value = 309.3125
$ws.Range("J131").Value = 941.6548
$ws.Range("K131").Value = 927.9375
$ws.Range("L131").Value = 2824.9644
$ws.Range("M131").Value = 4112.0625
$ws.Range("N131").Value = -12904.9644
$ws.Range("H132").Value = 50001016
$ws.Range("I132").Value = 90910296
$ws.Range("J132").Value = 787.6667
$ws.Range("K132").Value = 818192664
$ws.Range("L132").Value = 7089.0003
$ws.Range("M132").Value = -818190134
$ws.Range("N132").Value = -12149.0003

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 174572.28
$ws.Range("I80").Value = 5500
$ws.Range("J80").Value = 400002
$ws.Range("K80").Value = 5500
$ws.Range("L80").Value = 400002
$ws.Range("M80").Value = -4502
$ws.Range("N80").Value = -401998
$ws.Range("H83").Value = 174572.28
$ws.Range("I83").Value = 5500
$ws.Range("J83").Value = 400002
$ws.Range("K83").Value = 27500
$ws.Range("L83").Value = 2000010
$ws.Range("M83").Value = -22508
$ws.Range("N83").Value = -2009994
$ws.Range("H132").Value = 15153428
$ws.Range("I132").Value = 34484656
$ws.Range("J132").Value = 1923.6757
$ws.Range("K132").Value = 103453968
$ws.Range("L132").Value = 5771.0271
$ws.Range("M132").Value = -103451438
$ws.Range("N132").Value = -10831.0271

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H125").Value = 47418.75
$ws.Range("J125").Value = 47418.75
$ws.Range("L125").Value = 47418.75
$ws.Range("N125").Value = -57258.75
$ws.Range("H132").Value = 4389.8887
$ws.Range("I132").Value = 4848.614
$ws.Range("J132").Value = 3669.0356
$ws.Range("K132").Value = 14545.842
$ws.Range("L132").Value = 11007.1068
$ws.Range("M132").Value = -12015.842
$ws.Range("N132").Value = -16067.1068
